# Update cryptos worksheet with the latest scraped market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.094.34'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').Value = '1.792.67'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = "'221.82"
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = "'0.549"
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = "'31.50"
$ws.Range('E8').Value = '  -4.79%  '
$ws.Range('D9').Value = "'0.289"
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').Value = "'0.0710"
$ws.Range('E10').Value = '  +6.88%  '
$ws.Range('D11').Value = "'0.0925"
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').Value = '2.065.54'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = '1.805.29'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').Value = "'10.63"
$ws.Range('E14').Value = '  -4.53%  '
$ws.Range('D15').Value = "'0.627"
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '34.128.20'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = "'4.23"
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = "'68.40"
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').Value = "'245.07"
$ws.Range('E19').Value = '  -4.21%  '
$ws.Range('D20').Value = '0.0₃0784'
$ws.Range('E20').Value = '  +4.66%  '
$ws.Range('D21').Value = "'0.997"
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').Value = "'10.71"
$ws.Range('E22').Value = '  +2.31%  '
$ws.Range('D23').Value = "'4.14"
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').Value = "'2.13"
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = "'158.38"
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').Value = "'16.34"
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('D27').Value = "'7.04"
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').Value = "'0.0522"
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').Value = "'3.71"
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('D32').Value = "'1.20"
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('D33').Value = "'3.52"
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').Value = "'1.84"
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('D35').Value = '1.405.26'
$ws.Range('E35').Value = '  -3.78%  '
$ws.Range('D36').Value = "'1.05"
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('D37').Value = "'0.628"
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('D38').Value = "'0.0186"
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').Value = "'0.936"
$ws.Range('E39').Value = '  +4.23%  '
$ws.Range('D40').Value = "'79.77"
$ws.Range('E40').Value = '  -3.90%  '
$ws.Range('D41').Value = "'2.72"
$ws.Range('E41').Value = '  -4.84%  '
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').Value = "'2.12"
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = "'0.0497"
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'5.88"
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('D47').Value = '1.961.34'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('D48').Value = "'105.68"
$ws.Range('E48').Value = '  +5.71%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').Value = "'11.69"
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  +5.86%  '
